# Auto-generated edit script: refreshes cached market-price figures
# (columns H-N) across all 8 job sheets per the scheduled-runner update.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1837.359
$ws.Range("I15").Value = 1837.359
$ws.Range("K15").Value = 5512.076999999999
$ws.Range("M15").Value = -5343.076999999999
$ws.Range("H33").Value = 874.7692
$ws.Range("I33").Value = 486.44446
$ws.Range("J33").Value = 1748.5
$ws.Range("K33").Value = 486.44446
$ws.Range("L33").Value = 1748.5
$ws.Range("M33").Value = -257.44446
$ws.Range("N33").Value = -2206.5
$ws.Range("H98").Value = 1562.125
$ws.Range("I98").Value = 1500
$ws.Range("J98").Value = 1997
$ws.Range("K98").Value = 1500
$ws.Range("L98").Value = 1997
$ws.Range("M98").Value = -2
$ws.Range("N98").Value = -4993
$ws.Range("H116").Value = 17265.334
$ws.Range("J116").Value = 6665.3335
$ws.Range("L116").Value = 6665.3335
$ws.Range("N116").Value = -13549.3335
$ws.Range("H122").Value = 1562.125
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 1997
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 5991
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -10891
$ws.Range("H137").Value = 4632687.5
$ws.Range("I137").Value = 2962
$ws.Range("J137").Value = 13892138
$ws.Range("K137").Value = 8886
$ws.Range("L137").Value = 41676414
$ws.Range("M137").Value = -6336
$ws.Range("N137").Value = -41681514

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3228.9375
$ws.Range("I2").Value = 2886.4
$ws.Range("K2").Value = 2886.4
$ws.Range("M2").Value = -2773.4
$ws.Range("H32").Value = 28760856
$ws.Range("I32").Value = 29114178
$ws.Range("K32").Value = 29114178
$ws.Range("M32").Value = -29113891
$ws.Range("H45").Value = 4887.8237
$ws.Range("I45").Value = 4661.154
$ws.Range("K45").Value = 4661.154
$ws.Range("M45").Value = -4284.154
$ws.Range("H74").Value = 3149.3333
$ws.Range("I74").Value = 3361.4285
$ws.Range("J74").Value = 2407
$ws.Range("K74").Value = 3361.4285
$ws.Range("L74").Value = 2407
$ws.Range("M74").Value = -2487.4285
$ws.Range("N74").Value = -4155
$ws.Range("H77").Value = 3149.3333
$ws.Range("I77").Value = 3361.4285
$ws.Range("J77").Value = 2407
$ws.Range("K77").Value = 16807.1425
$ws.Range("L77").Value = 12035
$ws.Range("M77").Value = -12439.1425
$ws.Range("N77").Value = -20771
$ws.Range("H97").Value = 3500
$ws.Range("I97").Value = 3000
$ws.Range("J97").Value = 3750
$ws.Range("K97").Value = 3000
$ws.Range("L97").Value = 3750
$ws.Range("M97").Value = -2504
$ws.Range("N97").Value = -4742
$ws.Range("H110").Value = 1836.8
$ws.Range("I110").Value = 1295
$ws.Range("J110").Value = 2649.5
$ws.Range("K110").Value = 1295
$ws.Range("L110").Value = 2649.5
$ws.Range("M110").Value = 750
$ws.Range("N110").Value = -6739.5
$ws.Range("H116").Value = 3228.9375
$ws.Range("I116").Value = 2886.4
$ws.Range("K116").Value = 2886.4
$ws.Range("M116").Value = -592.4000000000001
$ws.Range("H122").Value = 2463
$ws.Range("I122").Value = 1868.5
$ws.Range("K122").Value = 5605.5
$ws.Range("M122").Value = -3155.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3228.9375
$ws.Range("I3").Value = 2886.4
$ws.Range("K3").Value = 2886.4
$ws.Range("M3").Value = -2772.4
$ws.Range("H86").Value = 7835.6665
$ws.Range("I86").Value = 8750
$ws.Range("J86").Value = 6007
$ws.Range("K86").Value = 8750
$ws.Range("L86").Value = 6007
$ws.Range("M86").Value = -7627
$ws.Range("N86").Value = -8253
$ws.Range("H89").Value = 7835.6665
$ws.Range("I89").Value = 8750
$ws.Range("J89").Value = 6007
$ws.Range("K89").Value = 43750
$ws.Range("L89").Value = 30035
$ws.Range("M89").Value = -38134
$ws.Range("N89").Value = -41267
$ws.Range("H105").Value = 2901.5
$ws.Range("I105").Value = 1773.3334
$ws.Range("K105").Value = 1773.3334
$ws.Range("M105").Value = -26.33339999999998
$ws.Range("I107").Value = 731.7273
$ws.Range("J107").Value = 1991.3334
$ws.Range("K107").Value = 731.7273
$ws.Range("L107").Value = 1991.3334
$ws.Range("M107").Value = 1188.2727
$ws.Range("N107").Value = -5831.3334

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5909.607
$ws.Range("I31").Value = 1989
$ws.Range("J31").Value = 6761.913
$ws.Range("K31").Value = 1989
$ws.Range("L31").Value = 6761.913
$ws.Range("M31").Value = -1694
$ws.Range("N31").Value = -7351.913
$ws.Range("H34").Value = 5909.607
$ws.Range("I34").Value = 1989
$ws.Range("J34").Value = 6761.913
$ws.Range("K34").Value = 1989
$ws.Range("L34").Value = 6761.913
$ws.Range("M34").Value = -1787
$ws.Range("N34").Value = -7165.913
$ws.Range("H107").Value = 1392.15
$ws.Range("I107").Value = 612.2
$ws.Range("K107").Value = 612.2
$ws.Range("M107").Value = 1307.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 3579.4211
$ws.Range("J109").Value = 4307.6924
$ws.Range("L109").Value = 12923.0772
$ws.Range("N109").Value = -15003.0772
$ws.Range("H132").Value = 436150.8
$ws.Range("I132").Value = 691.1667
$ws.Range("K132").Value = 6220.5003
$ws.Range("M132").Value = -3690.5003

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 40000
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 40000
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 40000
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -40926
$ws.Range("H48").Value = 30000
$ws.Range("J48").Value = 30000
$ws.Range("L48").Value = 30000
$ws.Range("N48").Value = -30970
$ws.Range("H49").Value = 30000
$ws.Range("J49").Value = 30000
$ws.Range("L49").Value = 30000
$ws.Range("N49").Value = -30368
$ws.Range("H97").Value = 1301.3077
$ws.Range("I97").Value = 810.7273
$ws.Range("J97").Value = 3999.5
$ws.Range("K97").Value = 810.7273
$ws.Range("L97").Value = 3999.5
$ws.Range("M97").Value = -314.7273
$ws.Range("N97").Value = -4991.5
$ws.Range("H102").Value = 2308.5
$ws.Range("I102").Value = 2120.5557
$ws.Range("K102").Value = 2120.5557
$ws.Range("M102").Value = -498.5556999999999
$ws.Range("H122").Value = 998
$ws.Range("I122").Value = 998
$ws.Range("K122").Value = 2994
$ws.Range("M122").Value = -544

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 30033
$ws.Range("I41").Value = 30033
$ws.Range("K41").Value = 30033
$ws.Range("M41").Value = -29595
$ws.Range("H55").Value = 497
$ws.Range("J55").Value = 517.6667
$ws.Range("L55").Value = 517.6667
$ws.Range("N55").Value = -863.6667
$ws.Range("H82").Value = 2561.111
$ws.Range("I82").Value = 2690
$ws.Range("K82").Value = 2690
$ws.Range("M82").Value = -2329
$ws.Range("H85").Value = 2561.111
$ws.Range("I85").Value = 2690
$ws.Range("K85").Value = 2690
$ws.Range("M85").Value = -1442
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 73295.336
$ws.Range("J112").Value = 73295.336
$ws.Range("L112").Value = 73295.336
$ws.Range("N112").Value = -76249.336
$ws.Range("H132").Value = 1905.7727
$ws.Range("I132").Value = 1651.1875
$ws.Range("K132").Value = 4953.5625
$ws.Range("M132").Value = -2423.5625

